# Commit: "Add data tsv for cafe and minor changes"
#
# - cafe (row 7, samples_retained) gets its real sample counts: D7 ("#neg.")
#   was a placeholder formula (=936-72-144-144); replace with the literal
#   value 720 now that the tsv data is in. G7 ("n", a shared SUM formula)
#   recalculates on its own once D7 changes.
# - I7 held a "TODO" note (highlighted yellow) reminding the author notes
#   were still needed; the note is now resolved, so clear the cell's
#   content and remove the yellow highlight, leaving the (now unused)
#   style in place.
# - Leave the selection on G7, matching the author's last cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("samples_retained")

# D7: replace the placeholder formula with the real literal value.
$ws.Range("D7").Value = 720

# I7: the "TODO" note is resolved -- clear the text but keep the cell's
# style slot, then turn off the yellow highlight fill.
$ws.Range("I7").ClearContents()
$ws.Range("I7").Interior.Pattern = -4142   # xlPatternNone ("No Fill")

# Leave the cursor where the author left it.
$ws.Range("G7").Select()
